$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "2024-07-20T13:34:00.000Z"
$ws.Range("D5").Value = "2024-07-20T13:34:00.000Z"
$ws.Range("D6").Value = "2024-07-20T13:34:00.000Z"
$ws.Range("D8").Value = "2024-07-20T13:34:00.000Z"
$ws.Range("D12").Value = "2024-07-20T13:34:00.000Z"
$ws.Range("D13").Value = "2024-07-20T13:34:00.000Z"

$ws.Range("W13").Value = 63522000
$ws.Range("AA13").Value = 126758000
$ws.Range("AE13").Value = 190280000
$ws.Range("AH13").Value = 155580000
$ws.Range("AK13").Value = 28
$ws.Range("AQ13").Value = 181380000
